$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.958.88'
$ws.Range("E2").Value = '  +3.35%  '
$ws.Range("D3").Value = '2.445.69'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '321.09'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.61%  '
$ws.Range("E7").Value = '  +0.79%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.98%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.89'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.83%  '
$ws.Range("E11").Value = '  +0.64%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.123'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.88%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.17'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.69%  '
$ws.Range("E14").Value = '  +1.58%  '
$ws.Range("D15").Value = '2.828.09'
$ws.Range("E15").Value = '  +0.83%  '
$ws.Range("D16").Value = '2.447.80'
$ws.Range("E16").Value = '  -0.21%  '
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '45.795.67'
$ws.Range("E18").Value = '  +3.20%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.42'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '0.0₃0933'
$ws.Range("E21").Value = '  +2.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '71.24'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.37'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '246.64'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.44%  '
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.87'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("E28").Value = '  +0.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.68'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.66%  '
$ws.Range("E30").Value = '  +1.07%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '49.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.61%  '
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.95'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.34'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0758'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.53'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.15%  '
$ws.Range("E38").Value = '  -0.53%  '
$ws.Range("E39").Value = '  +0.93%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.52'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.60%  '
$ws.Range("E41").Value = '  +4.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.111'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.86'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.78%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0292'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.04%  '
$ws.Range("D45").Value = '1.957.20'
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("E46").Value = '  -2.93%  '
$ws.Range("E47").Value = '  +0.55%  '
$ws.Range("E48").Value = '  +11.05%  '
$ws.Range("E49").Value = '  -5.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.98'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.45%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.24'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.03%  '
